# Fruta / hortaliza, semanal
# Insert one new weekly record at row 224 (pushing the existing
# rows 224-348 down to 225-349, dimension grows from T348 to T349).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 224; this shifts all
# rows 224..348 down to 225..349 and preserves their contents.
$ws.Rows(224).Insert()

# Populate the newly inserted row 224 with the new record's data.
$ws.Range("A224").Value = 6
$ws.Range("B224").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C224").Value = "Metropolitana"
$ws.Range("D224").Value = 45029
$ws.Range("E224").Value = 13
$ws.Range("F224").Value = "Fruta"
$ws.Range("G224").Value = 100101
$ws.Range("H224").Value = "Berries"
$ws.Range("I224").Value = 100101004
$ws.Range("J224").Value = "Frambuesa"
$ws.Range("K224").Value = "Sin especificar"
$ws.Range("L224").Value = "Primera"
$ws.Range("M224").Value = 150
$ws.Range("N224").Value = 8000
$ws.Range("O224").Value = 8000
$ws.Range("P224").Value = 8000
$ws.Range("Q224").Value = "`$/bandeja 2 kilos"
$ws.Range("R224").Value = "Provincia de Curicó"
$ws.Range("S224").Value = 4000
$ws.Range("T224").Value = 2
